# TC03_Canine_Filter_SamplePatho-Lymphoma.xlsx
# "Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study"
#
# The CasesTab Cypher query (cell B2 on the "startup" sheet) dropped its
# trailing OPTIONAL MATCH cohort column - the last RETURN line
#   coalesce(co.cohort_description, '') AS `Cohort`
# (and the now-trailing comma on the prior line) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nMATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [`"Lymphoma`"]  `nOPTIONAL MATCH (co:cohort)<-[*]-(c)`n  WITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $casesQuery

# Matches the saved workbook's new selection (was B4/A4-topLeft before the edit).
$ws.Range("B2").Select() | Out-Null
